$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B16").Value = "veh"
$ws.Range("B17").Value = "givelic"
$ws.Range("B18").Value = "newrentveh"

$ws.Range("D17").Value = "Выдача гос лицензий"
$ws.Range("D18").Value = "Создания арендуемого автомобиля"
$ws.Range("D16").Value = "Заспамить транспорт"

$ws.Range("D21").Select()
